# Generate Report for Handback
# Updates the handoff/handback timestamps and priority values produced by
# a fresh handback report generation run.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest HO Xliff Generate Date" column (G) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-23 18:16:22"
$wsOverview.Range("G5").Value = "2016-08-23 18:16:22"

# --- zh-cn sheet: Priority (E), Correspond Handoff Datetime (H),
#     Correspond Handback DateTime (K) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"
$wsZhCn.Range("H2").Value = "2016-08-23 18:16:17"
$wsZhCn.Range("H5").Value = "2016-08-23 18:16:17"
$wsZhCn.Range("K2").Value = "2016-08-23 18:16:35"
$wsZhCn.Range("K5").Value = "2016-08-23 18:16:35"

# --- de-de sheet: Priority (E), Correspond Handoff Datetime (H),
#     Correspond Handback DateTime (K) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"
$wsDeDe.Range("H2").Value = "2016-08-23 18:16:22"
$wsDeDe.Range("H5").Value = "2016-08-23 18:16:22"
$wsDeDe.Range("K2").Value = "2016-08-23 18:16:42"
$wsDeDe.Range("K5").Value = "2016-08-23 18:16:42"
